# "Generate Report for Handback" - refresh the localization-status report
# after a successful handback: the handback files are now in sync with
# en-US, so the "Ready for handoff" status becomes "Handed back: in sync
# with en-US", the Latest Handback DateTime stamps move forward, and the
# stale "handback file is not the latest" error clears out.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("K2").Value = "2016-08-24 12:52:22"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).AutoFit()
$zhcn.Columns.Item(16).AutoFit()

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("K2").Value = "2016-08-24 12:52:28"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).AutoFit()
$dede.Columns.Item(16).AutoFit()
